# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Price cells that look like plain decimals are forced to Text format first
# so Excel doesn't "helpfully" reinterpret them as numbers (which would
# drop trailing zeros / introduce floating point noise) - e.g. "604.40"
# must stay the literal string "604.40", not the number 604.4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.065.32'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').Value = '3.585.04'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.40'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '195.77'
$ws.Range('E6').Value = '  -1.70%  '
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -3.40%  '
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('E11').Value = '  -1.21%  '
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.55'
$ws.Range('D14').Value = '4.146.89'
$ws.Range('E14').Value = '  +2.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '598.87'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.04'
$ws.Range('E16').Value = '  +2.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.28'
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('D18').Value = '70.235.26'
$ws.Range('E18').Value = '  +0.39%  '
$ws.Range('D19').Value = '3.587.08'
$ws.Range('E19').Value = '  +2.25%  '
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.994'
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.88'
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '102.75'
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.13'
$ws.Range('E24').Value = '  +0.65%  '
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.09'
$ws.Range('E26').Value = '  -1.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.84'
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.59'
$ws.Range('E28').Value = '  -2.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.73'
$ws.Range('E29').Value = '  -0.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.12'
$ws.Range('E30').Value = '  -1.88%  '
$ws.Range('E31').Value = '  -6.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.37'
$ws.Range('E32').Value = '  -3.12%  '
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.48'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('D35').Value = '3.891.22'
$ws.Range('E35').Value = '  +4.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.24'
$ws.Range('E36').Value = '  +7.00%  '
$ws.Range('E37').Value = '  +2.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '526.41'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.393'
$ws.Range('E40').Value = '  +0.19%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.61'
$ws.Range('E41').Value = '  +1.11%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.04'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('E43').Value = '  -2.23%  '
$ws.Range('E44').Value = '  -1.90%  '
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.33'
$ws.Range('E47').Value = '  +0.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.58'
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000249'
$ws.Range('E50').Value = '  +3.73%  '
